$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 250
$ws1.Range("F3").Value = 77
$ws1.Range("F4").Value = 9514
$ws1.Range("F5").Value = 619
$ws1.Range("F6").Value = 162
$ws1.Range("F7").Value = 293
$ws1.Range("F8").Value = 358
$ws1.Range("F9").Value = 412
$ws1.Range("F11").Value = 193
$ws1.Range("F13").Value = 448
$ws1.Range("F14").Value = 12154
$ws1.Range("F19").Value = 241
$ws1.Range("F20").Value = 40
$ws1.Range("F21").Value = 175
$ws1.Range("F24").Value = 2729
$ws1.Range("F25").Value = 2100
$ws1.Range("F26").Value = 71
$ws1.Range("F28").Value = 54
$ws1.Range("F29").Value = 2147
$ws1.Range("F30").Value = 1012
$ws1.Range("F31").Value = 4204
$ws1.Range("F32").Value = 3657
$ws1.Range("F33").Value = 575
$ws1.Range("F35").Value = 3061
$ws1.Range("F36").Value = 30
$ws1.Range("F37").Value = 1326
$ws1.Range("F39").Value = 774
$ws1.Range("F40").Value = 22
$ws1.Range("F41").Value = 110
$ws1.Range("F42").Value = 426
$ws1.Range("F43").Value = 542
$ws1.Range("F44").Value = 70
$ws1.Range("F45").Value = 138
$ws1.Range("F46").Value = 230
$ws1.Range("F48").Value = 129
$ws1.Range("F49").Value = 143

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F17").Value = 26

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 250
$ws4.Range("F5").Value = 77
$ws4.Range("F6").Value = 9514
$ws4.Range("F7").Value = 619
$ws4.Range("F9").Value = 162
$ws4.Range("F10").Value = 293
$ws4.Range("F11").Value = 358
$ws4.Range("F12").Value = 412
$ws4.Range("F14").Value = 193
$ws4.Range("F15").Value = 448
$ws4.Range("F16").Value = 12154
$ws4.Range("F20").Value = 241
$ws4.Range("F22").Value = 175
$ws4.Range("F25").Value = 2729
$ws4.Range("F26").Value = 2100
$ws4.Range("F27").Value = 71
$ws4.Range("F29").Value = 54
$ws4.Range("F30").Value = 2147
$ws4.Range("F31").Value = 1012
$ws4.Range("F32").Value = 4204
$ws4.Range("F33").Value = 3657
$ws4.Range("F34").Value = 575
$ws4.Range("F36").Value = 3061
$ws4.Range("F37").Value = 1326
$ws4.Range("F39").Value = 774
$ws4.Range("F40").Value = 110
$ws4.Range("F41").Value = 426
$ws4.Range("F43").Value = 542
$ws4.Range("F44").Value = 70
$ws4.Range("F45").Value = 138
$ws4.Range("F46").Value = 230
$ws4.Range("F48").Value = 129
$ws4.Range("F49").Value = 143
